# Update the "StatQuery" cell for the ParticipantsTab row (C2) with the new
# SQL text: the Files count now combines sequencing + pathology files, an
# extra "sex_at_birth = 'Female'" filter was added, and the query is
# rewritten to drive FROM df_study (left-joined out) instead of
# df_participant.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newStatQuery = @"
SELECT
    COUNT(DISTINCT std.study_ID) AS "Studies",
    COUNT(DISTINCT prt.participant_id) AS "Participants",
    COUNT(DISTINCT smp.sample_id) AS "Samples",
    (COUNT(DISTINCT seq.id) + COUNT(DISTINCT paf.id)) AS "Files"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN 
    df_sequencing_file seq ON smp.id = seq."sample.id"
LEFT JOIN 
    df_pathology_file paf ON smp.id = paf."sample.id"
WHERE 
    std.study_ID = 'phs002430' 
    AND prt.race = 'Asian' 
    AND prt.sex_at_birth = 'Female';
"@

$ws.Range("C2").Value = $newStatQuery

# The author's workbook was left scrolled/selected at C2 (the cell they just
# edited) instead of the original C6.
$ws.Range("C2").Select()
